$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New delivery rows appended to the log (rows 29-35), matching the
# "Date Delivered to Xin" export for 2024-07-01 (serial 45474).
# Values are written in the specific column order below so that the
# shared-string table is built up in the same sequence as the source
# workbook (people first, then sensor codes, then raw device ids).

# --- Column G (person who delivered it) ---
$ws.Cells.Item(29,7).Value = "Stewart Norrish "
$ws.Cells.Item(31,7).Value = "Stewart  Norrish "
$ws.Cells.Item(32,7).Value = "Ella Stankiewicz "
$ws.Cells.Item(30,7).Value = "Mark Bjornestad "
$ws.Cells.Item(35,7).Value = "Ella "
$ws.Cells.Item(33,7).Value = "Stewart  Norrish "
$ws.Cells.Item(34,7).Value = "Mark Bjornestad "

# --- Column A (device / sensor name) ---
$ws.Cells.Item(29,1).Value = "WS17-WLE"
$ws.Cells.Item(31,1).Value = "PT7-VD9"
$ws.Cells.Item(32,1).Value = "WS48-ABS"
$ws.Cells.Item(33,1).Value = "WS27-XLU"
$ws.Cells.Item(30,1).Value = "WS10-C67"

# --- Columns B (rh id) and C (t id) for rows 29-33 ---
$ws.Cells.Item(29,2).Value = "65c29a0b39f940000d986a32"
$ws.Cells.Item(29,3).Value = "65c29a0ba6f677000ddd0937"
$ws.Cells.Item(30,2).Value = "667dd9e00cb3d2000e8458f0"
$ws.Cells.Item(30,3).Value = "667dd9e0dea5cc000edfc671"
$ws.Cells.Item(31,2).Value = "648b4aa004a601000d19a6f7"
$ws.Cells.Item(31,3).Value = "648b4aa166b6ec000ff468ec"
$ws.Cells.Item(32,2).Value = "6644f0dd573ffb000ce2b97f"
$ws.Cells.Item(32,3).Value = "6644f0ddb921b4000b90b85d"
$ws.Cells.Item(33,2).Value = "65d6457c7a715d000bf94dc0"
$ws.Cells.Item(33,3).Value = "65d6457d7a715d000c7d068c"

# --- Rows 34-35: device name + ids together ---
$ws.Cells.Item(34,1).Value = "WS37-4NE"
$ws.Cells.Item(34,2).Value = "6601cdd45f4c803643ea5668"
$ws.Cells.Item(34,3).Value = "6601cdd45fe5e134ed5d28be"
$ws.Cells.Item(35,1).Value = "WS12-ESF"
$ws.Cells.Item(35,2).Value = "6633cda1efd480018020fc75"
$ws.Cells.Item(35,3).Value = "6633cda2972270000bb4ffca"

# --- Latitude / longitude / date columns for all new rows ---
$ws.Cells.Item(29,4).Value = 40.316295459999999
$ws.Cells.Item(29,5).Value = -103.5642296
$ws.Cells.Item(29,6).Value = 45474
$ws.Cells.Item(29,6).NumberFormat = "d-mmm-yy"

$ws.Cells.Item(30,4).Value = 44.695286350000003
$ws.Cells.Item(30,5).Value = -108.7679644
$ws.Cells.Item(30,6).Value = 45474
$ws.Cells.Item(30,6).NumberFormat = "d-mmm-yy"

$ws.Cells.Item(31,4).Value = 40.144148600000001
$ws.Cells.Item(31,5).Value = -102.6879031
$ws.Cells.Item(31,6).Value = 45474
$ws.Cells.Item(31,6).NumberFormat = "d-mmm-yy"

$ws.Cells.Item(32,4).Value = 40.290162100000003
$ws.Cells.Item(32,5).Value = -104.8992951
$ws.Cells.Item(32,6).Value = 45474
$ws.Cells.Item(32,6).NumberFormat = "d-mmm-yy"

$ws.Cells.Item(33,4).Value = 39.962213920000003
$ws.Cells.Item(33,5).Value = -102.29797979999999
$ws.Cells.Item(33,6).Value = 45474
$ws.Cells.Item(33,6).NumberFormat = "d-mmm-yy"

$ws.Cells.Item(34,4).Value = 44.5044781
$ws.Cells.Item(34,5).Value = -108.3984542
$ws.Cells.Item(34,6).Value = 45474
$ws.Cells.Item(34,6).NumberFormat = "d-mmm-yy"

$ws.Cells.Item(35,4).Value = 40.135103569999998
$ws.Cells.Item(35,5).Value = -105.0301903
$ws.Cells.Item(35,6).Value = 45474
$ws.Cells.Item(35,6).NumberFormat = "d-mmm-yy"

# Scroll the view down to the newly added rows and select the first new
# entry, mirroring the author's on-screen state when they saved.
$ws.Range("A28:E28").Select()
